$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four observation rows (17-20) had their record identities (Id / Ost / Nord
# and a couple of optional attribute cells) reshuffled between rows while every
# other column (species, locality, dates, observer, ...) stayed put:
#   after-row 17 <- before-row 18
#   after-row 18 <- before-row 20
#   after-row 19 <- before-row 17   (brings the Enhet/Alder-Stadium/Kon/Metod/
#                                     Bestamningsmetod cells along with it)
#   after-row 20 <- before-row 19

# Row 17: becomes the old row 18's identity
$ws.Range("A17").Value = 111821926
$ws.Range("Q17").Value = 550846.2444635418
$ws.Range("R17").Value = 6681625.195240833
# Row 17 no longer carries the optional J/K/L/N/AF attribute cells
$ws.Range("J17").ClearContents()
$ws.Range("K17").ClearContents()
$ws.Range("L17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("AF17").ClearContents()

# Row 18: becomes the old row 20's identity
$ws.Range("A18").Value = 111821927
$ws.Range("Q18").Value = 550819.8901872271
$ws.Range("R18").Value = 6681733.007140613

# Row 19: becomes the old row 17's identity, including its extra attribute cells
$ws.Range("A19").Value = 111821924
$ws.Range("Q19").Value = 550675.3931295178
$ws.Range("R19").Value = 6681937.422269406
$ws.Range("J19").ClearContents()
$ws.Range("K19").Value = "blomning"
$ws.Range("L19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("AF19").ClearContents()

# Row 20: becomes the old row 19's identity
$ws.Range("A20").Value = 111821928
$ws.Range("Q20").Value = 550825.9503372401
$ws.Range("R20").Value = 6681726.144349095
